# Titman (2015): incorporate bookings/substitutions/score information into the
# tied_teams team orderings used downstream for the suspense/surprise/shock
# probability calculations. Column O on Sheet1 is "tied_teams".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40-52: ['Costa Rica', 'Ireland'] -> ['Ireland', 'Costa Rica']
$ws.Range("O40:O52").Value = "['Ireland', 'Costa Rica']"

# Rows 53-59: ['Costa Rica', 'Colombia', 'Ireland', 'Argentina'] -> ['Ireland', 'Colombia', 'Costa Rica', 'Argentina']
$ws.Range("O53:O59").Value = "['Ireland', 'Colombia', 'Costa Rica', 'Argentina']"

# Rows 60-62: ['Colombia', 'Argentina'] is unchanged.

# Rows 63-73: ['Scotland', 'Austria', 'Colombia', 'Argentina'] -> ['Colombia', 'Scotland', 'Argentina', 'Austria']
$ws.Range("O63:O73").Value = "['Colombia', 'Scotland', 'Argentina', 'Austria']"

# Row 78: ['South Korea', 'Netherlands'] -> ['Netherlands', 'South Korea']
$ws.Range("O78").Value = "['Netherlands', 'South Korea']"

# Rows 101-102: ['Netherlands', 'Italy'] is unchanged.

# Row 104: ['United States', 'Netherlands'] -> ['Netherlands', 'United States']
$ws.Range("O104").Value = "['Netherlands', 'United States']"
